$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list update (GitHub Actions data refresh)
# For Price (D) column: force text number-format before assignment so
# numeric-looking strings (e.g. "115.50", "1.00") are stored as text,
# matching the original inlineStr cells, then restore default style.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '51.322.65'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.82%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.756.71'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +3.38%  '

$ws.Range('E4').Value = '  -0.03%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '115.50'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.69%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '330.64'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.39%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.532'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.59%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.574'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.51%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.47'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.03%  '

$ws.Range('E11').Value = '  +1.13%  '

$ws.Range('E12').Value = '  +0.64%  '

$ws.Range('E13').Value = '  +3.10%  '

$ws.Range('E14').Value = '  +4.18%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.187.30'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.20%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.750.73'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.24%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.889'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.58%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '51.292.51'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.89%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.62'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.76%  '

$ws.Range('E20').Value = '  +4.75%  '

$ws.Range('E21').Value = '  +1.25%  '

$ws.Range('E22').Value = '  +0.63%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '280.19'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.96%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '70.22'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.03%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.61'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.31%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '27.02'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.70%  '

$ws.Range('E27').Value = '  +0.12%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.36'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.60%  '

$ws.Range('E29').Value = '  -0.60%  '

$ws.Range('E30').Value = '  -1.10%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '35.72'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.28%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '50.07'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.40%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.63'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.75%  '

$ws.Range('E34').Value = '  +0.96%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '19.40'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.20%  '

$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.12'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.08%  '

$ws.Range('B37').Value = 'FirstDigitalUSD'
$ws.Range('C37').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.15%  '

$ws.Range('E38').Value = '  -1.14%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.24'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.49%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '129.46'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.74%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '23.79'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.76%  '

$ws.Range('E42').Value = '  +10.57%  '

$ws.Range('E43').Value = '  +3.84%  '

$ws.Range('E44').Value = '  +0.45%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.42'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.45%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.118.18'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.16%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.22'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +9.82%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.27'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.08%  '

$ws.Range('E49').Value = '  +3.21%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.08'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.49%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.53'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +8.81%  '
